$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.815.65"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "2.084.35"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'235.06"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'59.31"
$ws.Range("E7").Value = "  +3.62%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.391"
$ws.Range("D10").Value = "'0.0793"
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("D11").Value = "'0.106"
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("D12").Value = "2.393.64"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "'14.64"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "'21.46"
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("D15").Value = "'0.770"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "'5.32"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "2.086.58"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "37.782.96"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "'6.24"
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("D20").Value = "'71.81"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "'228.85"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'2.42"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "'170.85"
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("D27").Value = "'0.140"
$ws.Range("E27").Value = "  +9.21%  "
$ws.Range("D28").Value = "'9.04"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").Value = "'1.44"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'19.52"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "'4.69"
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("D33").Value = "'0.0631"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "'4.71"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").Value = "'2.53"
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("D36").Value = "'3.53"
$ws.Range("E36").Value = "  +7.39%  "
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'5.45"
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("D40").Value = "'0.0986"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "'99.10"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'2.94"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "'0.0215"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "1.468.47"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "'1.18"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'4.24"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").Value = "'16.08"
$ws.Range("E47").Value = "  +5.04%  "
$ws.Range("D48").Value = "'1.07"
$ws.Range("E48").Value = "  +3.81%  "
$ws.Range("D49").Value = "'7.47"
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").Value = "'48.19"
$ws.Range("E51").Value = "  +7.99%  "
